$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-12-31 12:54:29"

for ($r = 2; $r -le 397; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

$ws.Cells.Item(336, 13).Value = "Pasquier Pancakes Choco 10 Stück - Online kein Bestand 4.30 Schweizer Franken"
